$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.464.97"
$ws.Range("E2").Value = "  -5.51%  "

# Row 3
$ws.Range("D3").Value = "3.352.37"
$ws.Range("E3").Value = "  -6.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'557.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.76%  "

# Row 6
$ws.Range("D6").Value = "'182.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.82%  "

# Row 7
$ws.Range("D7").Value = "'0.595"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.10%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "3.342.76"
$ws.Range("E9").Value = "  -6.53%  "

# Row 10
$ws.Range("E10").Value = "  -13.76%  "

# Row 11
$ws.Range("D11").Value = "'0.590"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.33%  "

# Row 12
$ws.Range("D12").Value = "'47.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.50%  "

# Row 13
$ws.Range("E13").Value = "  -10.92%  "

# Row 14
$ws.Range("D14").Value = "'8.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.79%  "

# Row 15
$ws.Range("D15").Value = "3.886.32"
$ws.Range("E15").Value = "  -6.60%  "

# Row 16
$ws.Range("D16").Value = "'595.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -14.52%  "

# Row 17
$ws.Range("D17").Value = "66.260.55"
$ws.Range("E17").Value = "  -5.90%  "

# Row 18
$ws.Range("D18").Value = "3.360.16"
$ws.Range("E18").Value = "  -6.75%  "

# Row 19
$ws.Range("E19").Value = "  -4.52%  "

# Row 20
$ws.Range("E20").Value = "  -7.06%  "

# Row 21
$ws.Range("D21").Value = "'11.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.86%  "

# Row 22
$ws.Range("D22").Value = "'0.907"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.57%  "

# Row 23
$ws.Range("D23").Value = "'16.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.43%  "

# Row 24
$ws.Range("D24").Value = "'5.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.24%  "

# Row 25
$ws.Range("D25").Value = "'96.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.17%  "

# Row 26
$ws.Range("E26").Value = "  -10.04%  "

# Row 27
$ws.Range("D27").Value = "'2.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.78%  "

# Row 28
$ws.Range("D28").Value = "'9.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.41%  "

# Row 29
$ws.Range("D29").Value = "'8.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.14%  "

# Row 30
$ws.Range("D30").Value = "'30.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.54%  "

# Row 31
$ws.Range("D31").Value = "'3.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.55%  "

# Row 32
$ws.Range("E32").Value = "  -10.72%  "

# Row 33
$ws.Range("E33").Value = "  -9.07%  "

# Row 34
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "3.830.77"
$ws.Range("E34").Value = "  +0.77%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.105"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.70%  "

# Row 36
$ws.Range("E36").Value = "  -8.39%  "

# Row 37
$ws.Range("D37").Value = "'530.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.92%  "

# Row 38
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("D39").Value = "'3.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +33.50%  "

# Row 40
$ws.Range("E40").Value = "  -7.59%  "

# Row 41
$ws.Range("E41").Value = "  -14.63%  "

# Row 42
$ws.Range("E42").Value = "  -6.91%  "

# Row 43
$ws.Range("E43").Value = "  -10.22%  "

# Row 44
$ws.Range("E44").Value = "  -8.70%  "

# Row 45
$ws.Range("D45").Value = "'32.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.19%  "

# Row 46
$ws.Range("D46").Value = "'0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.44%  "

# Row 47
$ws.Range("E47").Value = "  -12.57%  "

# Row 48
$ws.Range("E48").Value = "  -8.84%  "

# Row 49
$ws.Range("D49").Value = "'0.129"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.22%  "

# Row 50
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

# Row 51
$ws.Range("D51").Value = "'7.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.16%  "
